# Deploying to gh-pages from  @ f61b6b856fe916e5aba0a2d26945892c8efe815c
# Updates the 16.5.1.1a indicator sheet:
#  - refreshes the title strings in A1/B1/C1 (Kyrgyz/Russian/English)
#  - adds the 2020 data column (I)
#  - updates the active cell selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title row: normalize "16.5.1.1a." -> "16.5.1.1a " and refresh wording ---
$ws.Range("A1").Value = "16.5.1.1a ""Аткаруу бийлигинин мамлекеттик органдарындагы жана жергиликтүү өз алдынча башкаруу органдарындагы коррупциянын деңгээли жөнүндө жеке түшүнүк"" индекси"
$ws.Range("B1").Value = "16.5.1.1a Индекс ""Личное представление об уровне коррупции в государственных органах исполнительной власти и органах местного самоуправления''"
$ws.Range("C1").Value = "16.5.1.1a Index ""Personal views about the level of corruption in executive government authorities and local government''"

# --- Add 2020 column (I) mirroring the formatting of the 2019 column (H) ---
$ws.Range("H4").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("I4").Value = 2020

$ws.Range("H5").Copy()
$ws.Range("I5").PasteSpecial(-4122)
$ws.Range("I5").Value = 12.3
$ws.Range("I5").NumberFormat = "0.0"

$ws.Range("H6").Copy()
$ws.Range("I6").PasteSpecial(-4122)
$ws.Range("I6").Value = 40.3
$ws.Range("I6").NumberFormat = "0.0"

$ws.Range("H7").Copy()
$ws.Range("I7").PasteSpecial(-4122)
$ws.Range("I7").Value = 36.2
$ws.Range("I7").NumberFormat = "0.0"

$ws.Range("H8").Copy()
$ws.Range("I8").PasteSpecial(-4122)
$ws.Range("I8").Value = 44.3
$ws.Range("I8").NumberFormat = "0.0"

$ws.Range("H9").Copy()
$ws.Range("I9").PasteSpecial(-4122)
$ws.Range("I9").Value = 36
$ws.Range("I9").NumberFormat = "0.0"

$ws.Range("H10").Copy()
$ws.Range("I10").PasteSpecial(-4122)
$ws.Range("I10").Value = 2.7
$ws.Range("I10").NumberFormat = "0.0"

$ws.Range("H11").Copy()
$ws.Range("I11").PasteSpecial(-4122)
$ws.Range("I11").Value = 32.9
$ws.Range("I11").NumberFormat = "0.0"

$ws.Range("H12").Copy()
$ws.Range("I12").PasteSpecial(-4122)
$ws.Range("I12").Value = 11.3
$ws.Range("I12").NumberFormat = "0.0"

$ws.Range("H13").Copy()
$ws.Range("I13").PasteSpecial(-4122)
$ws.Range("I13").Value = -18.2
$ws.Range("I13").NumberFormat = "0.0"

$ws.Range("H14").Copy()
$ws.Range("I14").PasteSpecial(-4122)
$ws.Range("I14").Value = 33
$ws.Range("I14").NumberFormat = "0.0"

# --- Match the author's final UI selection ---
$ws.Range("F16").Select()
